$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Jamal Murray'
$ws.Range("B2").Value = 'PG,SG'
$ws.Range("C2").Value = 'Denver Nuggets'
$ws.Range("A3").Value = 'Kyrie Irving'
$ws.Range("B3").Value = 'PG,SG'
$ws.Range("C3").Value = 'Dallas Mavericks'
$ws.Range("A4").Value = 'CJ McCollum'
$ws.Range("B4").Value = 'PG,SG'
$ws.Range("C4").Value = 'New Orleans Pelicans'
$ws.Range("A5").Value = 'Kyle Kuzma'
$ws.Range("B5").Value = 'SF,PF'
$ws.Range("C5").Value = 'Milwaukee Bucks'
$ws.Range("A6").Value = 'Rui Hachimura'
$ws.Range("B6").Value = 'SF,PF'
$ws.Range("C6").Value = 'Los Angeles Lakers'
$ws.Range("A7").Value = 'Moses Moody'
$ws.Range("B7").Value = 'SG,SF'
$ws.Range("C7").Value = 'Golden State Warriors'
$ws.Range("A8").Value = 'Royce O''Neale'
$ws.Range("B8").Value = 'SF,PF'
$ws.Range("C8").Value = 'Phoenix Suns'
$ws.Range("A9").Value = 'Jordan Clarkson'
$ws.Range("B9").Value = 'SG,SF'
$ws.Range("C9").Value = 'Utah Jazz'
$ws.Range("A10").Value = 'Cason Wallace'
$ws.Range("B10").Value = 'PG,SG'
$ws.Range("C10").Value = 'Oklahoma City Thunder'
$ws.Range("A11").Value = 'Bam Adebayo'
$ws.Range("B11").Value = 'PF,C'
$ws.Range("C11").Value = 'Miami Heat'
$ws.Range("A12").Value = 'Shai Gilgeous-Alexander'
$ws.Range("B12").Value = 'PG,SG'
$ws.Range("C12").Value = 'Oklahoma City Thunder'
$ws.Range("A13").Value = 'Jordan Poole'
$ws.Range("B13").Value = 'PG,SG'
$ws.Range("C13").Value = 'Washington Wizards'
$ws.Range("A14").Value = 'Zach LaVine'
$ws.Range("B14").Value = 'SG,SF'
$ws.Range("C14").Value = 'Sacramento Kings'
$ws.Range("A15").Value = 'Tobias Harris'
$ws.Range("B15").Value = 'SF,PF'
$ws.Range("C15").Value = 'Detroit Pistons'
$ws.Range("A16").Value = 'John Collins'
$ws.Range("B16").Value = 'PF,C'
$ws.Range("C16").Value = 'Utah Jazz'
$ws.Range("A17").Value = 'Joel Embiid'
$ws.Range("B17").Value = 'C'
$ws.Range("C17").Value = 'Philadelphia 76ers'
$ws.Range("A18").Value = 'Lauri Markkanen'
$ws.Range("B18").Value = 'SF,PF'
$ws.Range("C18").Value = 'Utah Jazz'
